$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 277.9  # H9: 110.5 -> 277.9
$ws.Cells.Item(9, 9).Value = 111.28571  # I9: 69.14286 -> 111.28571
$ws.Cells.Item(9, 10).Value = 666.6667  # J9: 400 -> 666.6667
$ws.Cells.Item(9, 11).Value = 111.28571  # K9: 69.14286 -> 111.28571
$ws.Cells.Item(9, 12).Value = 666.6667  # L9: 400 -> 666.6667
$ws.Cells.Item(9, 13).Value = 57.71429000000001  # M9: 99.85714 -> 57.71429000000001
$ws.Cells.Item(9, 14).Value = -1004.6667  # N9: -738 -> -1004.6667

$ws.Cells.Item(11, 8).Value = 252.5  # H11: 183.33333 -> 252.5
$ws.Cells.Item(11, 9).Value = 252.5  # I11: 183.33333 -> 252.5
$ws.Cells.Item(11, 11).Value = 252.5  # K11: 183.33333 -> 252.5
$ws.Cells.Item(11, 13).Value = -112.5  # M11: -43.33332999999999 -> -112.5

$ws.Cells.Item(40, 8).Value = 1621.0769  # H40: 1694.6 -> 1621.0769
$ws.Cells.Item(40, 9).Value = 1008.2222  # I40: 1039.875 -> 1008.2222
$ws.Cells.Item(40, 10).Value = 3000  # J40: 2442.8572 -> 3000
$ws.Cells.Item(40, 11).Value = 1008.2222  # K40: 1039.875 -> 1008.2222
$ws.Cells.Item(40, 12).Value = 3000  # L40: 2442.8572 -> 3000
$ws.Cells.Item(40, 13).Value = -833.2222  # M40: -864.875 -> -833.2222
$ws.Cells.Item(40, 14).Value = -3350  # N40: -2792.8572 -> -3350

$ws.Cells.Item(43, 8).Value = 2399.6667  # H43: 2499.6 -> 2399.6667
$ws.Cells.Item(43, 10).Value = 2579.6  # J43: 2749.5 -> 2579.6
$ws.Cells.Item(43, 12).Value = 2579.6  # L43: 2749.5 -> 2579.6
$ws.Cells.Item(43, 14).Value = -2717.6  # N43: -2887.5 -> -2717.6

$ws.Cells.Item(55, 8).Value = 348.8889  # H55: 588.8 -> 348.8889
$ws.Cells.Item(55, 10).Value = 86.666664  # J55: 162 -> 86.666664
$ws.Cells.Item(55, 12).Value = 86.666664  # L55: 162 -> 86.666664
$ws.Cells.Item(55, 14).Value = -514.666664  # N55: -590 -> -514.666664

$ws.Cells.Item(69, 8).Value = 1523.9062  # H69: 1591.8518 -> 1523.9062
$ws.Cells.Item(69, 10).Value = 1498.871  # J69: 1564.6154 -> 1498.871
$ws.Cells.Item(69, 12).Value = 4496.613  # L69: 4693.8462 -> 4496.613
$ws.Cells.Item(69, 14).Value = -6244.613  # N69: -6441.8462 -> -6244.613

$ws.Cells.Item(72, 8).Value = 1523.9062  # H72: 1591.8518 -> 1523.9062
$ws.Cells.Item(72, 10).Value = 1498.871  # J72: 1564.6154 -> 1498.871
$ws.Cells.Item(72, 12).Value = 13489.839  # L72: 14081.5386 -> 13489.839
$ws.Cells.Item(72, 14).Value = -22225.839  # N72: -22817.5386 -> -22225.839

$ws.Cells.Item(76, 8).Value = 2780681.2  # H76: 3159.2778 -> 2780681.2
$ws.Cells.Item(76, 9).Value = 2881.4285  # I76: 3061.1667 -> 2881.4285
$ws.Cells.Item(76, 10).Value = 4276419.5  # J76: 3208.3333 -> 4276419.5
$ws.Cells.Item(76, 11).Value = 2881.4285  # K76: 3061.1667 -> 2881.4285
$ws.Cells.Item(76, 12).Value = 4276419.5  # L76: 3208.3333 -> 4276419.5
$ws.Cells.Item(76, 13).Value = -2566.4285  # M76: -2746.1667 -> -2566.4285
$ws.Cells.Item(76, 14).Value = -4277049.5  # N76: -3838.3333 -> -4277049.5

$ws.Cells.Item(79, 8).Value = 2780681.2  # H79: 3159.2778 -> 2780681.2
$ws.Cells.Item(79, 9).Value = 2881.4285  # I79: 3061.1667 -> 2881.4285
$ws.Cells.Item(79, 10).Value = 4276419.5  # J79: 3208.3333 -> 4276419.5
$ws.Cells.Item(79, 11).Value = 2881.4285  # K79: 3061.1667 -> 2881.4285
$ws.Cells.Item(79, 12).Value = 4276419.5  # L79: 3208.3333 -> 4276419.5
$ws.Cells.Item(79, 13).Value = -1789.4285  # M79: -1969.1667 -> -1789.4285
$ws.Cells.Item(79, 14).Value = -4278603.5  # N79: -5392.3333 -> -4278603.5

$ws.Cells.Item(107, 8).Value = 457.75  # H107: 469.2 -> 457.75
$ws.Cells.Item(107, 9).Value = 457.75  # I107: 469.2 -> 457.75
$ws.Cells.Item(107, 11).Value = 457.75  # K107: 469.2 -> 457.75
$ws.Cells.Item(107, 13).Value = 1462.25  # M107: 1450.8 -> 1462.25

$ws.Cells.Item(129, 8).Value = 1667358.4  # H129: 909880.6 -> 1667358.4
$ws.Cells.Item(129, 10).Value = 1667358.4  # J129: 909880.6 -> 1667358.4
$ws.Cells.Item(129, 12).Value = 5002075.199999999  # L129: 2729641.8 -> 5002075.199999999
$ws.Cells.Item(129, 14).Value = -5012075.199999999  # N129: -2739641.8 -> -5012075.199999999

$ws.Cells.Item(132, 8).Value = 4608.6113  # H132: 4452.8945 -> 4608.6113
$ws.Cells.Item(132, 9).Value = 4730.3335  # I132: 4537.8125 -> 4730.3335
$ws.Cells.Item(132, 11).Value = 14191.0005  # K132: 13613.4375 -> 14191.0005
$ws.Cells.Item(132, 13).Value = -11661.0005  # M132: -11083.4375 -> -11661.0005

$ws.Cells.Item(137, 8).Value = 49154.906  # H137: 64241.5 -> 49154.906
$ws.Cells.Item(137, 9).Value = 1329.5  # I137: 1558.125 -> 1329.5
$ws.Cells.Item(137, 10).Value = 112922.11  # J137: 126924.875 -> 112922.11
$ws.Cells.Item(137, 11).Value = 3988.5  # K137: 4674.375 -> 3988.5
$ws.Cells.Item(137, 12).Value = 338766.33  # L137: 380774.625 -> 338766.33
$ws.Cells.Item(137, 13).Value = -1438.5  # M137: -2124.375 -> -1438.5
$ws.Cells.Item(137, 14).Value = -343866.33  # N137: -385874.625 -> -343866.33

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(16, 8).Value = 825  # H16: 0 -> 825
$ws.Cells.Item(16, 9).Value = 650  # I16: 0 -> 650
$ws.Cells.Item(16, 10).Value = 1000  # J16: 0 -> 1000
$ws.Cells.Item(16, 11).Value = 650  # K16: 0 -> 650
$ws.Cells.Item(16, 12).Value = 1000  # L16: 0 -> 1000
$ws.Cells.Item(16, 13).Value = -363  # M16: None -> -363
$ws.Cells.Item(16, 14).Value = -1574  # N16: None -> -1574

$ws.Cells.Item(45, 8).Value = 4048.5925  # H45: 4158.9165 -> 4048.5925
$ws.Cells.Item(45, 9).Value = 4049.8333  # I45: 4260 -> 4049.8333
$ws.Cells.Item(45, 10).Value = 4047.6  # J45: 4086.7144 -> 4047.6
$ws.Cells.Item(45, 11).Value = 4049.8333  # K45: 4260 -> 4049.8333
$ws.Cells.Item(45, 12).Value = 4047.6  # L45: 4086.7144 -> 4047.6
$ws.Cells.Item(45, 13).Value = -3672.8333  # M45: -3883 -> -3672.8333
$ws.Cells.Item(45, 14).Value = -4801.6  # N45: -4840.7144 -> -4801.6

$ws.Cells.Item(74, 8).Value = 40002680  # H74: 50002810 -> 40002680
$ws.Cells.Item(74, 9).Value = 55558412  # I74: 62503196 -> 55558412
$ws.Cells.Item(74, 10).Value = 2230.5715  # J74: 1253.5 -> 2230.5715
$ws.Cells.Item(74, 11).Value = 55558412  # K74: 62503196 -> 55558412
$ws.Cells.Item(74, 12).Value = 2230.5715  # L74: 1253.5 -> 2230.5715
$ws.Cells.Item(74, 13).Value = -55557538  # M74: -62502322 -> -55557538
$ws.Cells.Item(74, 14).Value = -3978.5715  # N74: -3001.5 -> -3978.5715

$ws.Cells.Item(77, 8).Value = 40002680  # H77: 50002810 -> 40002680
$ws.Cells.Item(77, 9).Value = 55558412  # I77: 62503196 -> 55558412
$ws.Cells.Item(77, 10).Value = 2230.5715  # J77: 1253.5 -> 2230.5715
$ws.Cells.Item(77, 11).Value = 277792060  # K77: 312515980 -> 277792060
$ws.Cells.Item(77, 12).Value = 11152.8575  # L77: 6267.5 -> 11152.8575
$ws.Cells.Item(77, 13).Value = -277787692  # M77: -312511612 -> -277787692
$ws.Cells.Item(77, 14).Value = -19888.8575  # N77: -15003.5 -> -19888.8575

$ws.Cells.Item(97, 8).Value = 854.3333  # H97: 923.625 -> 854.3333
$ws.Cells.Item(97, 9).Value = 1102.2307  # I97: 1248.091 -> 1102.2307
$ws.Cells.Item(97, 11).Value = 1102.2307  # K97: 1248.091 -> 1102.2307
$ws.Cells.Item(97, 13).Value = -606.2307000000001  # M97: -752.0909999999999 -> -606.2307000000001

$ws.Cells.Item(132, 8).Value = 13942.098  # H132: 15790.777 -> 13942.098
$ws.Cells.Item(132, 9).Value = 1564.2858  # I132: 1767.0435 -> 1564.2858
$ws.Cells.Item(132, 11).Value = 4692.857400000001  # K132: 5301.1305 -> 4692.857400000001
$ws.Cells.Item(132, 13).Value = -2162.857400000001  # M132: -2771.1305 -> -2162.857400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 980.32654  # H94: 995.5417 -> 980.32654
$ws.Cells.Item(94, 9).Value = 955.15  # I94: 973.2308 -> 955.15
$ws.Cells.Item(94, 11).Value = 955.15  # K94: 973.2308 -> 955.15
$ws.Cells.Item(94, 13).Value = -504.15  # M94: -522.2308 -> -504.15

$ws.Cells.Item(99, 8).Value = 2082.6365  # H99: 1660.6 -> 2082.6365
$ws.Cells.Item(99, 9).Value = 1582  # I99: 1101.1111 -> 1582
$ws.Cells.Item(99, 11).Value = 1582  # K99: 1101.1111 -> 1582
$ws.Cells.Item(99, 13).Value = -84  # M99: 396.8888999999999 -> -84

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 19983.334  # H31: 21052.941 -> 19983.334
$ws.Cells.Item(31, 10).Value = 2971.2856  # J31: 3166.5 -> 2971.2856
$ws.Cells.Item(31, 12).Value = 2971.2856  # L31: 3166.5 -> 2971.2856
$ws.Cells.Item(31, 14).Value = -3561.2856  # N31: -3756.5 -> -3561.2856

$ws.Cells.Item(34, 8).Value = 19983.334  # H34: 21052.941 -> 19983.334
$ws.Cells.Item(34, 10).Value = 2971.2856  # J34: 3166.5 -> 2971.2856
$ws.Cells.Item(34, 12).Value = 2971.2856  # L34: 3166.5 -> 2971.2856
$ws.Cells.Item(34, 14).Value = -3375.2856  # N34: -3570.5 -> -3375.2856

$ws.Cells.Item(94, 8).Value = 9302.182000000001  # H94: 9306.637000000001 -> 9302.182000000001
$ws.Cells.Item(94, 9).Value = 2987.5  # I94: 2999.75 -> 2987.5
$ws.Cells.Item(94, 11).Value = 2987.5  # K94: 2999.75 -> 2987.5
$ws.Cells.Item(94, 13).Value = -2536.5  # M94: -2548.75 -> -2536.5

$ws.Cells.Item(105, 8).Value = 6579746.5  # H105: 6945275 -> 6579746.5
$ws.Cells.Item(105, 10).Value = 1124  # J105: 1251.5714 -> 1124
$ws.Cells.Item(105, 12).Value = 1124  # L105: 1251.5714 -> 1124
$ws.Cells.Item(105, 14).Value = -4618  # N105: -4745.5714 -> -4618

$ws.Cells.Item(132, 8).Value = 17901.605  # H132: 18451.906 -> 17901.605
$ws.Cells.Item(132, 9).Value = 20434.852  # I132: 21209.576 -> 20434.852
$ws.Cells.Item(132, 11).Value = 61304.556  # K132: 63628.728 -> 61304.556
$ws.Cells.Item(132, 13).Value = -58774.556  # M132: -61098.728 -> -58774.556

$ws.Cells.Item(134, 8).Value = 1000.5263  # H134: 1084.1945 -> 1000.5263
$ws.Cells.Item(134, 9).Value = 886.1579  # I134: 1029.8334 -> 886.1579
$ws.Cells.Item(134, 10).Value = 1114.8948  # J134: 1138.5555 -> 1114.8948
$ws.Cells.Item(134, 11).Value = 2658.4737  # K134: 3089.5002 -> 2658.4737
$ws.Cells.Item(134, 12).Value = 3344.6844  # L134: 3415.6665 -> 3344.6844
$ws.Cells.Item(134, 13).Value = -123.4737  # M134: -554.5001999999999 -> -123.4737
$ws.Cells.Item(134, 14).Value = -8414.6844  # N134: -8485.666499999999 -> -8414.6844

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1200.875  # H5: 1329.5714 -> 1200.875
$ws.Cells.Item(5, 9).Value = 658.1429000000001  # I5: 717.8333 -> 658.1429000000001
$ws.Cells.Item(5, 11).Value = 1974.4287  # K5: 2153.4999 -> 1974.4287
$ws.Cells.Item(5, 13).Value = -1862.4287  # M5: -2041.4999 -> -1862.4287

$ws.Cells.Item(14, 8).Value = 279.08334  # H14: 280 -> 279.08334
$ws.Cells.Item(14, 9).Value = 279.08334  # I14: 280 -> 279.08334
$ws.Cells.Item(14, 11).Value = 837.2500200000001  # K14: 840 -> 837.2500200000001
$ws.Cells.Item(14, 13).Value = -664.2500200000001  # M14: -667 -> -664.2500200000001

$ws.Cells.Item(70, 8).Value = 3487.25  # H70: 4643.3 -> 3487.25
$ws.Cells.Item(70, 9).Value = 2571.2  # I70: 3079.111 -> 2571.2
$ws.Cells.Item(70, 10).Value = 5014  # J70: 5923.091 -> 5014
$ws.Cells.Item(70, 11).Value = 7713.599999999999  # K70: 9237.332999999999 -> 7713.599999999999
$ws.Cells.Item(70, 12).Value = 15042  # L70: 17769.273 -> 15042
$ws.Cells.Item(70, 13).Value = -7398.599999999999  # M70: -8922.332999999999 -> -7398.599999999999
$ws.Cells.Item(70, 14).Value = -15672  # N70: -18399.273 -> -15672

$ws.Cells.Item(73, 8).Value = 3487.25  # H73: 4643.3 -> 3487.25
$ws.Cells.Item(73, 9).Value = 2571.2  # I73: 3079.111 -> 2571.2
$ws.Cells.Item(73, 10).Value = 5014  # J73: 5923.091 -> 5014
$ws.Cells.Item(73, 11).Value = 7713.599999999999  # K73: 9237.332999999999 -> 7713.599999999999
$ws.Cells.Item(73, 12).Value = 15042  # L73: 17769.273 -> 15042
$ws.Cells.Item(73, 13).Value = -6621.599999999999  # M73: -8145.332999999999 -> -6621.599999999999
$ws.Cells.Item(73, 14).Value = -17226  # N73: -19953.273 -> -17226

$ws.Cells.Item(97, 8).Value = 1472  # H97: 1488.9231 -> 1472
$ws.Cells.Item(97, 9).Value = 690  # I97: 793.3333 -> 690
$ws.Cells.Item(97, 10).Value = 1614.1818  # J97: 1697.6 -> 1614.1818
$ws.Cells.Item(97, 11).Value = 2070  # K97: 2379.9999 -> 2070
$ws.Cells.Item(97, 12).Value = 4842.5454  # L97: 5092.799999999999 -> 4842.5454
$ws.Cells.Item(97, 13).Value = -1574  # M97: -1883.9999 -> -1574
$ws.Cells.Item(97, 14).Value = -5834.5454  # N97: -6084.799999999999 -> -5834.5454

$ws.Cells.Item(122, 8).Value = 721.1429000000001  # H122: 661.625 -> 721.1429000000001
$ws.Cells.Item(122, 10).Value = 1316  # J122: 1048.25 -> 1316
$ws.Cells.Item(122, 12).Value = 11844  # L122: 9434.25 -> 11844
$ws.Cells.Item(122, 14).Value = -16744  # N122: -14334.25 -> -16744

$ws.Cells.Item(131, 8).Value = 714.4299999999999  # H131: 720.1799999999999 -> 714.4299999999999
$ws.Cells.Item(131, 10).Value = 715.5859  # J131: 721.3939 -> 715.5859
$ws.Cells.Item(131, 12).Value = 2146.7577  # L131: 2164.1817 -> 2146.7577
$ws.Cells.Item(131, 14).Value = -12226.7577  # N131: -12244.1817 -> -12226.7577

$ws.Cells.Item(132, 8).Value = 1156.5834  # H132: 1460 -> 1156.5834
$ws.Cells.Item(132, 9).Value = 828.4286  # I132: 1200 -> 828.4286
$ws.Cells.Item(132, 11).Value = 7455.8574  # K132: 10800 -> 7455.8574
$ws.Cells.Item(132, 13).Value = -4925.8574  # M132: -8270 -> -4925.8574

$ws.Cells.Item(135, 8).Value = 1200.875  # H135: 1329.5714 -> 1200.875
$ws.Cells.Item(135, 9).Value = 658.1429000000001  # I135: 717.8333 -> 658.1429000000001
$ws.Cells.Item(135, 11).Value = 5923.2861  # K135: 6460.4997 -> 5923.2861
$ws.Cells.Item(135, 13).Value = -3388.2861  # M135: -3925.4997 -> -3388.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 2293.3333  # H14: 2297.1428 -> 2293.3333
$ws.Cells.Item(14, 10).Value = 2293.3333  # J14: 2297.1428 -> 2293.3333
$ws.Cells.Item(14, 12).Value = 2293.3333  # L14: 2297.1428 -> 2293.3333
$ws.Cells.Item(14, 14).Value = -2637.3333  # N14: -2641.1428 -> -2637.3333

$ws.Cells.Item(21, 8).Value = 1702.5  # H21: 0 -> 1702.5
$ws.Cells.Item(21, 10).Value = 1702.5  # J21: 0 -> 1702.5
$ws.Cells.Item(21, 12).Value = 1702.5  # L21: 0 -> 1702.5
$ws.Cells.Item(21, 14).Value = -2050.5  # N21: None -> -2050.5

$ws.Cells.Item(55, 8).Value = 149.28572  # H55: 73.789474 -> 149.28572
$ws.Cells.Item(55, 9).Value = 88.333336  # I55: 38.5 -> 88.333336
$ws.Cells.Item(55, 10).Value = 195  # J55: 113 -> 195
$ws.Cells.Item(55, 11).Value = 88.333336  # K55: 38.5 -> 88.333336
$ws.Cells.Item(55, 12).Value = 195  # L55: 113 -> 195
$ws.Cells.Item(55, 13).Value = 84.666664  # M55: 134.5 -> 84.666664
$ws.Cells.Item(55, 14).Value = -541  # N55: -459 -> -541

$ws.Cells.Item(61, 8).Value = 2714.9412  # H61: 2788.1516 -> 2714.9412
$ws.Cells.Item(61, 9).Value = 1148.1305  # I61: 1186.7273 -> 1148.1305
$ws.Cells.Item(61, 11).Value = 1148.1305  # K61: 1186.7273 -> 1148.1305
$ws.Cells.Item(61, 13).Value = -946.1305  # M61: -984.7273 -> -946.1305

$ws.Cells.Item(82, 8).Value = 2135.3684  # H82: 2236.2727 -> 2135.3684
$ws.Cells.Item(82, 9).Value = 2619.6667  # I82: 2400 -> 2619.6667
$ws.Cells.Item(82, 10).Value = 1699.5  # J82: 2099.8333 -> 1699.5
$ws.Cells.Item(82, 11).Value = 2619.6667  # K82: 2400 -> 2619.6667
$ws.Cells.Item(82, 12).Value = 1699.5  # L82: 2099.8333 -> 1699.5
$ws.Cells.Item(82, 13).Value = -2258.6667  # M82: -2039 -> -2258.6667
$ws.Cells.Item(82, 14).Value = -2421.5  # N82: -2821.8333 -> -2421.5

$ws.Cells.Item(85, 8).Value = 2135.3684  # H85: 2236.2727 -> 2135.3684
$ws.Cells.Item(85, 9).Value = 2619.6667  # I85: 2400 -> 2619.6667
$ws.Cells.Item(85, 10).Value = 1699.5  # J85: 2099.8333 -> 1699.5
$ws.Cells.Item(85, 11).Value = 2619.6667  # K85: 2400 -> 2619.6667
$ws.Cells.Item(85, 12).Value = 1699.5  # L85: 2099.8333 -> 1699.5
$ws.Cells.Item(85, 13).Value = -1371.6667  # M85: -1152 -> -1371.6667
$ws.Cells.Item(85, 14).Value = -4195.5  # N85: -4595.8333 -> -4195.5

$ws.Cells.Item(113, 8).Value = 2714.9412  # H113: 2788.1516 -> 2714.9412
$ws.Cells.Item(113, 9).Value = 1148.1305  # I113: 1186.7273 -> 1148.1305
$ws.Cells.Item(113, 11).Value = 1148.1305  # K113: 1186.7273 -> 1148.1305
$ws.Cells.Item(113, 13).Value = 1021.8695  # M113: 983.2727 -> 1021.8695

$ws.Cells.Item(132, 8).Value = 1826.2142  # H132: 1720.7742 -> 1826.2142
$ws.Cells.Item(132, 9).Value = 1245.0714  # I132: 1146.7222 -> 1245.0714
$ws.Cells.Item(132, 10).Value = 2407.3572  # J132: 2515.6155 -> 2407.3572
$ws.Cells.Item(132, 11).Value = 3735.2142  # K132: 3440.1666 -> 3735.2142
$ws.Cells.Item(132, 12).Value = 7222.071599999999  # L132: 7546.8465 -> 7222.071599999999
$ws.Cells.Item(132, 13).Value = -1205.2142  # M132: -910.1665999999996 -> -1205.2142
$ws.Cells.Item(132, 14).Value = -12282.0716  # N132: -12606.8465 -> -12282.0716

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 506  # H9: 0 -> 506
$ws.Cells.Item(9, 9).Value = 506  # I9: 0 -> 506
$ws.Cells.Item(9, 11).Value = 506  # K9: 0 -> 506
$ws.Cells.Item(9, 13).Value = -366  # M9: None -> -366

$ws.Cells.Item(14, 8).Value = 0  # H14: 4492 -> 0
$ws.Cells.Item(14, 9).Value = 0  # I14: 4004 -> 0
$ws.Cells.Item(14, 10).Value = 0  # J14: 4980 -> 0
$ws.Cells.Item(14, 11).Value = 0  # K14: 4004 -> 0
$ws.Cells.Item(14, 12).Value = 0  # L14: 4980 -> 0
$ws.Cells.Item(14, 13).Value = $null  # M14: -3836 -> (removed)
$ws.Cells.Item(14, 14).Value = $null  # N14: -5316 -> (removed)

$ws.Cells.Item(126, 8).Value = 923.5  # H126: 1273.36 -> 923.5
$ws.Cells.Item(126, 9).Value = 789.1667  # I126: 1408.6 -> 789.1667
$ws.Cells.Item(126, 10).Value = 1125  # J126: 1070.5 -> 1125
$ws.Cells.Item(126, 11).Value = 2367.5001  # K126: 4225.799999999999 -> 2367.5001
$ws.Cells.Item(126, 12).Value = 3375  # L126: 3211.5 -> 3375
$ws.Cells.Item(126, 13).Value = 102.4998999999998  # M126: -1755.799999999999 -> 102.4998999999998
$ws.Cells.Item(126, 14).Value = -8315  # N126: -8151.5 -> -8315
